$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ask" sheet has two side-by-side English/French question-word tables:
#   A/B = first table (English | French), C/D = second table (English | French)
# This edit swaps each pair so French comes first: A<->B and C<->D.

# --- Swap columns A and B (both already have explicit custom widths) ---
$ws.Columns("A").Cut()
$ws.Columns("C").Insert()

# --- Swap columns C and D ---
# (Column D has no explicit width yet, so after the cut/insert shuffle the
# engine leaves a stray "non-custom" width marker on column C; clear it and
# restore the original centered cell style that ClearFormats wipes.)
$ws.Columns("D").Cut()
$ws.Columns("C").Insert()
$ws.Columns("C").ClearFormats()
$ws.Range("D1:D4").Copy()
$ws.Range("C1:C4").PasteSpecial(-4122)

# --- Match the author's final selection: whole column C ---
$ws.Columns("C").Select()
